$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 7
$ws.Range("G7").Value = 2
